$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.374.14'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.587.62'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.06'
$ws.Range('E5').Value = '  -3.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.95'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = '3.581.85'
$ws.Range('E7').Value = '  -2.02%  '
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.675'
$ws.Range('E10').Value = '  -3.44%  '
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.57'
$ws.Range('E12').Value = '  -3.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '4.152.94'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '3.599.91'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '67.322.12'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.30'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.39'
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('E21').Value = '  -3.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '401.97'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.47'
$ws.Range('E23').Value = '  +22.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.20'
$ws.Range('E24').Value = '  -3.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.61'
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.55'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.09'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.79'
$ws.Range('E29').Value = '  +3.73%  '
$ws.Range('E30').Value = '  +10.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.14'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.29'
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '669.00'
$ws.Range('E33').Value = '  +10.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.20'
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('E35').Value = '  +0.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.05'
$ws.Range('E36').Value = '  -2.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.49'
$ws.Range('E37').Value = '  -4.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.424'
$ws.Range('E38').Value = '  +8.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '0.0₃0785'
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.92'
$ws.Range('E41').Value = '  +15.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.11'
$ws.Range('E42').Value = '  +8.29%  '
$ws.Range('D43').Value = '3.159.41'
$ws.Range('E43').Value = '  +14.13%  '
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0418'
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.13'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.132'
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.76'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.36'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('E51').Value = '  -3.13%  '
